$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.722.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.209.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.41'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.39%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.542.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.214.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.781'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.621.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.62%  '
$ws.Range("E19").Value = '  +1.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E23").Value = '  +9.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.96%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("E28").Value = '  +21.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '170.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0799'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0325'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("E41").Value = '  +5.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '58.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.54%  '
$ws.Range("E45").Value = '  +2.54%  '
$ws.Range("E46").Value = '  +20.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0978'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.68%  '
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("E51").Value = '  +1.98%  '
